$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added on top of the existing data for this
# market/product. Insert a fresh row at position 7 (this shifts the
# previously-existing rows 7..101 down to 8..102, dimension grows to T102)
# and then populate the new row with the new record's values.
$ws.Rows("7:7").Insert()

$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44643
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100109
$ws.Cells.Item(7, 8).Value = "Uva"
$ws.Cells.Item(7, 9).Value = 100109001
$ws.Cells.Item(7, 10).Value = "Uva"
$ws.Cells.Item(7, 11).Value = "Red Globe"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 360
$ws.Cells.Item(7, 14).Value = 9000
$ws.Cells.Item(7, 15).Value = 10000
$ws.Cells.Item(7, 16).Value = 9500
$ws.Cells.Item(7, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(7, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(7, 19).Value = 528
$ws.Cells.Item(7, 20).Value = 18
